# Apply data corrections to the "sharp" (G) and "subpixel" (H) columns,
# and update the saved view/selection state of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rows where both G (sharp) and H (subpixel) go from 7 -> 8
$rowsSevenToEight = @(6, 25, 30, 31, 33, 40, 41, 43, 54, 55, 57, 58)
foreach ($r in $rowsSevenToEight) {
    $ws.Range("G$r").Value = 8
    $ws.Range("H$r").Value = 8
}

# Rows where only G (sharp) goes from 3 -> 5 (H left untouched)
$rowsThreeToFive = @(13, 24, 37)
foreach ($r in $rowsThreeToFive) {
    $ws.Range("G$r").Value = 5
}

# Update the sheet's view: clear the frozen/scrolled topLeftCell and move
# the active selection to I24.
$ws.Activate()
$ws.Range("I24").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
